# Fruta / hortaliza, semanal
# Re-sync the three data rows (2,3,4) so they are ordered by date:
#   Row2 <- old Row3 values, Row3 <- old Row4 values, Row4 <- old Row2 values
# Columns A,B,C,E,F,G,H,I,J,K,L,R are identical across the three rows already,
# so only D (Fecha), M (Volumen), N/O/P (Precios), Q (Unidad), S (Precio $/Kg)
# and T (Kg/unidad) need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> becomes what row 3 used to hold
$ws.Range("D2").Value = 44875
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("S2").Value = 1600
$ws.Range("T2").Value = 10

# Row 3 -> becomes what row 4 used to hold
$ws.Range("D3").Value = 44874
$ws.Range("M3").Value = 67
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1600
$ws.Range("T3").Value = 10

# Row 4 -> becomes what row 2 used to hold
$ws.Range("D4").Value = 44855
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/bandeja 5 kilos"
$ws.Range("S4").Value = 3000
$ws.Range("T4").Value = 5
